$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  2  = 5
  3  = 5
  4  = 3
  5  = 2
  6  = 1
  7  = 3
  8  = 7
  9  = 5
  10 = 9
  11 = 2
  12 = 3
  13 = 3
  14 = 3
  15 = 6
  16 = 1
  17 = 4
  18 = 2
  19 = 3
  20 = 9
  21 = 4
  22 = 6
  23 = 2
  24 = 4
  25 = 0
  26 = 7
  27 = 4
  28 = 3
  29 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
